$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'47.352.95"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "'2.492.35"
$ws.Range("E3").Value = "  -0.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'321.28"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("D6").Value = "'108.15"
$ws.Range("E6").Value = "  +2.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -1.06%  "

# Row 10
$ws.Range("D10").Value = "'38.74"
$ws.Range("E10").Value = "  +4.85%  "

# Row 11
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
$ws.Range("D12").Value = "'0.124"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13
$ws.Range("D13").Value = "'18.39"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14
$ws.Range("E14").Value = "  -1.49%  "

# Row 15
$ws.Range("D15").Value = "'2.881.26"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").Value = "'2.492.07"
$ws.Range("E16").Value = "  -1.04%  "

# Row 17
$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("D18").Value = "'47.261.04"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("E19").Value = "  +0.84%  "

# Row 20
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "  +0.50%  "

# Row 21
$ws.Range("D21").Value = "'0.0₃0932"
$ws.Range("E21").Value = "  -0.71%  "

# Row 22
$ws.Range("D22").Value = "'2.70"
$ws.Range("E22").Value = "  +13.77%  "

# Row 23
$ws.Range("D23").Value = "'70.33"
$ws.Range("E23").Value = "  -0.69%  "

# Row 24
$ws.Range("D24").Value = "'245.50"
$ws.Range("E24").Value = "  -2.51%  "

# Row 25
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +0.58%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").Value = "'25.74"
$ws.Range("E27").Value = "  -2.03%  "

# Row 28
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +3.97%  "

# Row 29
$ws.Range("E29").Value = "  -0.52%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  +1.93%  "

# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'34.59"
$ws.Range("E31").Value = "  -1.85%  "

# Row 32
$ws.Range("D32").Value = "'49.64"
$ws.Range("E32").Value = "  +0.16%  "

# Row 33
$ws.Range("D33").Value = "'20.82"
$ws.Range("E33").Value = "  +6.17%  "

# Row 34
$ws.Range("D34").Value = "'5.33"
$ws.Range("E34").Value = "  +0.55%  "

# Row 35
$ws.Range("E35").Value = "  +0.72%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").Value = "'4.67"
$ws.Range("E38").Value = "  +1.14%  "

# Row 39
$ws.Range("E39").Value = "  -1.35%  "

# Row 40
$ws.Range("D40").Value = "'23.24"
$ws.Range("E40").Value = "  +8.78%  "

# Row 41
$ws.Range("E41").Value = "  -0.41%  "

# Row 42
$ws.Range("E42").Value = "  -0.44%  "

# Row 43
$ws.Range("D43").Value = "'117.85"
$ws.Range("E43").Value = "  -3.26%  "

# Row 44
$ws.Range("D44").Value = "'0.0297"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45
$ws.Range("D45").Value = "'1.987.16"
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  +1.06%  "

# Row 47
$ws.Range("E47").Value = "  -6.04%  "

# Row 48
$ws.Range("D48").Value = "'9.06"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("E49").Value = "  -1.14%  "

# Row 50
$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = "  -5.75%  "

# Row 51
$ws.Range("D51").Value = "'56.77"
$ws.Range("E51").Value = "  +3.75%  "
